# Update odds values in the "Jogos da Semana" sheet for row 5, 6 and 7
# as described by the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 35
$ws.Range("H5").Value = 9.25
$ws.Range("J5").Value = 23
$ws.Range("K5").Value = 3.95
$ws.Range("L5").Value = 1.2
$ws.Range("P5").Value = 11.5
$ws.Range("Q5").Value = 1.12
$ws.Range("R5").Value = 4.4
$ws.Range("S5").Value = 1.09
$ws.Range("T5").Value = 6.1
$ws.Range("U5").Value = 2.54
$ws.Range("V5").Value = 1.49
$ws.Range("W5").Value = 200
$ws.Range("Y5").Value = 200
$ws.Range("AB5").Value = 500
$ws.Range("AD5").Value = 32
$ws.Range("AE5").Value = 65
$ws.Range("AI5").Value = 7.6
$ws.Range("AN5").Value = 40
$ws.Range("AO5").Value = 400
$ws.Range("AU5").Value = 15
$ws.Range("BA5").Value = 24

# Row 6 updates
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.92

# Row 7 updates
$ws.Range("J7").Value = 5.2
$ws.Range("AG7").Value = 300
